$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G2 : "Correspond Handoff Datetime" timestamp
$wsOverview.Range("G2").Value = "2016-08-23 03:04:18"

# zh-cn!H2 : "Correspond Handoff Datetime" timestamp
$wsZhCn.Range("H2").Value = "2016-08-23 03:04:13"

# zh-cn!K2 : "Correspond Handback DateTime" timestamp
$wsZhCn.Range("K2").Value = "2016-08-23 03:04:30"

# de-de!H2 : "Correspond Handoff Datetime" timestamp (same original text as Overview!G2)
$wsDeDe.Range("H2").Value = "2016-08-23 03:04:18"

# de-de!K2 : "Correspond Handback DateTime" timestamp
$wsDeDe.Range("K2").Value = "2016-08-23 03:04:37"
